$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Helper: write a value as literal TEXT (matches existing t="str" cells in this
# sheet) without creating a new cell style. We build it as a text formula
# (="...") then immediately collapse it to a static value via copy/paste-values,
# which leaves the cell's style index untouched (unlike NumberFormat="@").
function Set-TextValue($cell, $text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# New rows to append (A, C, F) - all values are stored as text in this sheet
$newRows = @(
    @{ Row = 62; A = $null; C = "397_红绿灵草_undefined_undefined_1bunch"; F = "10" },
    @{ Row = 63; A = $null; C = "490_米花 粉_rice flower pink_undefined_1bunch"; F = "8" },
    @{ Row = 64; A = $null; C = "44_拉丝粉_Spider Pink_Gerbera L._20stems"; F = "10" },
    @{ Row = 65; A = $null; C = "412_紫罗兰粉_violet pink_undefined_1bunch"; F = "25" },
    @{ Row = 66; A = $null; C = "411_紫罗兰白_violet white_undefined_1bunch"; F = "25" },
    @{ Row = 67; A = "20";  C = "535_雪果白_snow berry white_undefined_1bunch"; F = "5" },
    @{ Row = 68; A = $null; C = "411_紫罗兰白_violet white_undefined_1bunch"; F = "10" },
    @{ Row = 69; A = $null; C = "454_蓝星花_tweedia blue_undefined_1bunch"; F = "10" }
)

foreach ($r in $newRows) {
    if ($null -ne $r.A) {
        Set-TextValue $ws.Cells.Item($r.Row, 1) $r.A
    }
    Set-TextValue $ws.Cells.Item($r.Row, 3) $r.C
    Set-TextValue $ws.Cells.Item($r.Row, 6) $r.F
}

# Update the Summary sheet's G2 concatenated digit-string (same text value
# technique, since it is a 118-digit number stored as text).
$summary = $wb.Worksheets.Item("Summary")
$g2 = $summary.Cells.Item(2, 7)
Set-TextValue $g2 "0151961819423202321152252414103081177663240404011560542085565623555122054012420104025156035205583315354310810252551010"

# Used range grew from A1:L61 to A1:L69 - extend the "numbers stored as text"
# ignore-error flag (xlNumberAsText = 3) to cover the newly added rows too.
$usedRng = $ws.Range("A1:L69")
$usedRng.Errors.Item(3).Ignore = $true
